# Rename the "...RI" sheets to "...-RI" (add missing hyphen before RI suffix)
$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("TGZ-S-48-50_100RI").Name   = "TGZ-S-48-50_100-RI"
$wb.Worksheets.Item("TGZ-S-48-100_250RI").Name  = "TGZ-S-48-100_250-RI"
$wb.Worksheets.Item("TGZ-S-48-100_300RI").Name  = "TGZ-S-48-100_300-RI"

# Update the saved selection on the "TGZ-S-48-50_100-RI" sheet to G26
$ws = $wb.Worksheets.Item("TGZ-S-48-50_100-RI")
$ws.Select()
$ws.Range("G26").Select()
